# IPL Error Codes workbook update
# 1) Add a new "TOURNAMENT" sheet (error codes 740 series) after STAT.
# 2) Add "720 series" marker row to STAT and "640 series" marker row to TEAM.

$wb = $excel.ActiveWorkbook

# --- Step 1: STAT sheet gets a "720 series" note in A9 --------------------
$statSheet = $wb.Worksheets.Item("STAT")
$statSheet.Range("A9").Value = "720 series"

# --- Step 2: TEAM sheet gets a "640 series" note in A6 --------------------
$teamSheet = $wb.Worksheets.Item("TEAM")
$teamSheet.Range("A6").Value = "640 series"
$teamSheet.Range("A7").Select() | Out-Null

# --- Step 3: add the new TOURNAMENT sheet at the end of the workbook ------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "TOURNAMENT"

# Header row - copy formatting (bold/border/centered) from STAT's header
$statSheet.Range("A1:B1").Copy() | Out-Null
$newSheet.Range("A1:B1").PasteSpecial(-4122) | Out-Null
$newSheet.Range("A1").Value = "Err Code"
$newSheet.Range("B1").Value = "Error Description"

# Data row - copy formatting from MATCH's data row (uses the same non-centered
# bordered style as the other single-row sheets)
$matchSheet = $wb.Worksheets.Item("MATCH")
$matchSheet.Range("A2:B2").Copy() | Out-Null
$newSheet.Range("A2:B2").PasteSpecial(-4122) | Out-Null
$newSheet.Range("A2").Value = 741
$newSheet.Range("B2").Value = "Invalid tournament name"

# Column widths roughly matching the bestfit widths used elsewhere
$newSheet.Range("A1").ColumnWidth = 10
$newSheet.Range("B1").ColumnWidth = 23
$newSheet.Range("C1").ColumnWidth = 21.5

$newSheet.Range("C6").Select() | Out-Null

$excel.CutCopyMode = $false

Write-Output "Updated workbook: added TOURNAMENT sheet, STAT/TEAM series markers."
